$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Membre1")
$ws.Select()

# Copy formatting from row 6 (last filled log entry) down to the new row 7
$ws.Range("A6:D6").Copy()
$ws.Range("A7:D7").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in the new log entry
$ws.Range("A7").Value = 43031
$ws.Range("B7").Value = "30min"
$ws.Range("C7").Value = "Création du contrôle/model, ajout de la navigation entre les menus"
$ws.Range("D7").Value = "La première fois que j'utilise un dictionnaire, le find_if et des expressions lambda."

# Reset scroll position to the top and move the selection to the new row
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D7").Select()
